$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Monthly summary table (A2:B16): fill in the just-completed month's total steps ---
$ws.Cells.Item(10, 2).Value2 = 377377

# --- Weekly/daily log (A23:B53): roll the tracked window forward one month ---
# Newly logged daily steps, Jun 30 - Jul 5
$ws.Cells.Item(23, 1).Value2 = 42552
$ws.Cells.Item(23, 2).Value2 = 22431

$ws.Cells.Item(24, 1).Value2 = 42553
$ws.Cells.Item(24, 2).Value2 = 6609

$ws.Cells.Item(25, 1).Value2 = 42554
$ws.Cells.Item(25, 2).Value2 = 3391

$ws.Cells.Item(26, 1).Value2 = 42555
$ws.Cells.Item(26, 2).Value2 = 16090

$ws.Cells.Item(27, 1).Value2 = 42556
$ws.Cells.Item(27, 2).Value2 = 9756

$ws.Cells.Item(28, 1).Value2 = 42557
$ws.Cells.Item(28, 2).Value2 = 3204

# Remaining days of the new month: dates advance, no steps recorded yet so clear old values
for ($r = 29; $r -le 52; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r, 1).Value2 + 30
    $ws.Cells.Item($r, 2).ClearContents()
}

# New day row appended at the end of the (longer) month
$ws.Cells.Item(53, 1).Value2 = 42582
$ws.Cells.Item(53, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(53, 1).HorizontalAlignment = -4108

# --- View housekeeping ---
[void]$ws.Range("J30").Select()
